# Apply cell value updates to match the target revision.
# Values below were diffed from the canonical OOXML; each (sheet, cell) pair
# is set explicitly. A handful of cells are newly introduced (previously empty)
# and one cell (GSM!N98) is cleared back to empty.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 8567.333000000001
$ws.Range("I33").Value = 79.2
$ws.Range("J33").Value = 14630.286
$ws.Range("K33").Value = 79.2
$ws.Range("L33").Value = 14630.286
$ws.Range("M33").Value = 149.8
$ws.Range("N33").Value = -15088.286
$ws.Range("H105").Value = 45800
$ws.Range("J105").Value = 45800
$ws.Range("L105").Value = 45800
$ws.Range("N105").Value = -52788
$ws.Range("H112").Value = 2763.9268
$ws.Range("J112").Value = 2908.4736
$ws.Range("L112").Value = 8725.4208
$ws.Range("N112").Value = -10941.4208
$ws.Range("H129").Value = 981.44116
$ws.Range("I129").Value = 323.25
$ws.Range("J129").Value = 1069.2
$ws.Range("K129").Value = 969.75
$ws.Range("L129").Value = 3207.6
$ws.Range("M129").Value = 4030.25
$ws.Range("N129").Value = -13207.6
$ws.Range("H137").Value = 788.2195
$ws.Range("I137").Value = 692.1667
$ws.Range("J137").Value = 1479.8
$ws.Range("K137").Value = 2076.5001
$ws.Range("L137").Value = 4439.4
$ws.Range("M137").Value = 473.4998999999998
$ws.Range("N137").Value = -9539.4
$ws.Range("H138").Value = 3815.4878
$ws.Range("J138").Value = 5065.706
$ws.Range("L138").Value = 15197.118
$ws.Range("N138").Value = -25477.118

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 356457.25
$ws.Range("I32").Value = 3366.7693
$ws.Range("K32").Value = 3366.7693
$ws.Range("M32").Value = -3079.7693
$ws.Range("H45").Value = 2783.125
$ws.Range("I45").Value = 2674
$ws.Range("J45").Value = 2965
$ws.Range("K45").Value = 2674
$ws.Range("L45").Value = 2965
$ws.Range("M45").Value = -2297
$ws.Range("N45").Value = -3719
$ws.Range("H74").Value = 2471.4167
$ws.Range("I74").Value = 2594.1177
$ws.Range("K74").Value = 2594.1177
$ws.Range("M74").Value = -1720.1177
$ws.Range("H77").Value = 2471.4167
$ws.Range("I77").Value = 2594.1177
$ws.Range("K77").Value = 12970.5885
$ws.Range("M77").Value = -8602.588499999998
$ws.Range("H109").Value = 30770
$ws.Range("J109").Value = 30770
$ws.Range("L109").Value = 30770
$ws.Range("N109").Value = -33544
$ws.Range("H110").Value = 1292.8889
$ws.Range("I110").Value = 837.7692
$ws.Range("J110").Value = 2476.2
$ws.Range("K110").Value = 837.7692
$ws.Range("L110").Value = 2476.2
$ws.Range("M110").Value = 1207.2308
$ws.Range("N110").Value = -6566.2
$ws.Range("H132").Value = 1121.5555
$ws.Range("I132").Value = 773.7179599999999
$ws.Range("J132").Value = 3382.5
$ws.Range("K132").Value = 2321.15388
$ws.Range("L132").Value = 10147.5
$ws.Range("M132").Value = 208.8461200000002
$ws.Range("N132").Value = -15207.5
$ws.Range("H138").Value = 46644.145
$ws.Range("J138").Value = 46644.145
$ws.Range("L138").Value = 46644.145
$ws.Range("N138").Value = -56924.145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 232.10715
$ws.Range("I80").Value = 170.77777
$ws.Range("J80").Value = 261.1579
$ws.Range("K80").Value = 170.77777
$ws.Range("L80").Value = 261.1579
$ws.Range("M80").Value = 827.22223
$ws.Range("N80").Value = -2257.1579
$ws.Range("H83").Value = 232.10715
$ws.Range("I83").Value = 170.77777
$ws.Range("J83").Value = 261.1579
$ws.Range("K83").Value = 853.88885
$ws.Range("L83").Value = 1305.7895
$ws.Range("M83").Value = 4138.11115
$ws.Range("N83").Value = -11289.7895
$ws.Range("H86").Value = 3015.2896
$ws.Range("I86").Value = 2702.158
$ws.Range("J86").Value = 3328.4211
$ws.Range("K86").Value = 2702.158
$ws.Range("L86").Value = 3328.4211
$ws.Range("M86").Value = -1579.158
$ws.Range("N86").Value = -5574.4211
$ws.Range("H89").Value = 3015.2896
$ws.Range("I89").Value = 2702.158
$ws.Range("J89").Value = 3328.4211
$ws.Range("K89").Value = 13510.79
$ws.Range("L89").Value = 16642.1055
$ws.Range("M89").Value = -7894.789999999999
$ws.Range("N89").Value = -27874.1055

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1200
$ws.Range("J15").Value = 1200
$ws.Range("L15").Value = 1200
$ws.Range("N15").Value = -1540
$ws.Range("H58").Value = 1289.15
$ws.Range("I58").Value = 1047.2354
$ws.Range("K58").Value = 1047.2354
$ws.Range("M58").Value = -844.2354
$ws.Range("H122").Value = 1581.9615
$ws.Range("I122").Value = 1311.8422
$ws.Range("J122").Value = 2315.1428
$ws.Range("K122").Value = 3935.5266
$ws.Range("L122").Value = 6945.428400000001
$ws.Range("M122").Value = -1485.5266
$ws.Range("N122").Value = -11845.4284
$ws.Range("H132").Value = 1762.3684
$ws.Range("I132").Value = 1381.5883
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 4144.7649
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -1614.7649
$ws.Range("N132").Value = -20057
$ws.Range("H136").Value = 1289.15
$ws.Range("I136").Value = 1047.2354
$ws.Range("K136").Value = 3141.7062
$ws.Range("M136").Value = -591.7062000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 1857.1428
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 2000
$ws.Range("K20").Value = 3000
$ws.Range("L20").Value = 6000
$ws.Range("M20").Value = -2773
$ws.Range("N20").Value = -6454
$ws.Range("H56").Value = 4853412.5
$ws.Range("I56").Value = 4853412.5
$ws.Range("K56").Value = 4853412.5
$ws.Range("M56").Value = -4852882.5
$ws.Range("H74").Value = 6909.846
$ws.Range("I74").Value = 2503.25
$ws.Range("J74").Value = 8868.333000000001
$ws.Range("K74").Value = 7509.75
$ws.Range("L74").Value = 26604.999
$ws.Range("M74").Value = -6448.75
$ws.Range("N74").Value = -28726.999
$ws.Range("H77").Value = 6909.846
$ws.Range("I77").Value = 2503.25
$ws.Range("J77").Value = 8868.333000000001
$ws.Range("K77").Value = 22529.25
$ws.Range("L77").Value = 79814.997
$ws.Range("M77").Value = -17225.25
$ws.Range("N77").Value = -90422.997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 410
$ws.Range("I19").Value = 410
$ws.Range("K19").Value = 410
$ws.Range("M19").Value = -122
$ws.Range("H80").Value = 2856.6667
$ws.Range("I80").Value = 2022.2222
$ws.Range("J80").Value = 4108.3335
$ws.Range("K80").Value = 2022.2222
$ws.Range("L80").Value = 4108.3335
$ws.Range("M80").Value = -1024.2222
$ws.Range("N80").Value = -6104.3335
$ws.Range("H83").Value = 2856.6667
$ws.Range("I83").Value = 2022.2222
$ws.Range("J83").Value = 4108.3335
$ws.Range("K83").Value = 10111.111
$ws.Range("L83").Value = 20541.6675
$ws.Range("M83").Value = -5119.110999999999
$ws.Range("N83").Value = -30525.6675
$ws.Range("H92").Value = 9881.75
$ws.Range("J92").Value = 9881.75
$ws.Range("L92").Value = 9881.75
$ws.Range("N92").Value = -13625.75
$ws.Range("H93").Value = 20500
$ws.Range("J93").Value = 20500
$ws.Range("L93").Value = 20500
$ws.Range("N93").Value = -24244
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("N94").Value = -21352
$ws.Range("H98").Value = 0
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").Value = $null
$ws.Range("H99").Value = 7850
$ws.Range("I99").Value = 7850
$ws.Range("K99").Value = 7850
$ws.Range("M99").Value = -5604
$ws.Range("H101").Value = 50000
$ws.Range("J101").Value = 50000
$ws.Range("L101").Value = 50000
$ws.Range("N101").Value = -56490
$ws.Range("H141").Value = 46540
$ws.Range("J141").Value = 46540
$ws.Range("L141").Value = 46540
$ws.Range("N141").Value = -56900

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 576.0417
$ws.Range("I22").Value = 523.8889
$ws.Range("J22").Value = 732.5
$ws.Range("K22").Value = 523.8889
$ws.Range("L22").Value = 732.5
$ws.Range("M22").Value = -228.8889
$ws.Range("N22").Value = -1322.5
$ws.Range("H27").Value = 576.0417
$ws.Range("I27").Value = 523.8889
$ws.Range("J27").Value = 732.5
$ws.Range("K27").Value = 523.8889
$ws.Range("L27").Value = 732.5
$ws.Range("M27").Value = -416.8889
$ws.Range("N27").Value = -946.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H109").Value = 43028.5
$ws.Range("J109").Value = 43028.5
$ws.Range("L109").Value = 43028.5
$ws.Range("N109").Value = -45802.5
$ws.Range("H132").Value = 17546726
$ws.Range("I132").Value = 25642006
$ws.Range("J132").Value = 6953
$ws.Range("K132").Value = 76926018
$ws.Range("L132").Value = 20859
$ws.Range("M132").Value = -76923488
$ws.Range("N132").Value = -25919
